$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09476766666666665
$ws.Range("I2").Value = 0.02336090049363864
$ws.Range("J2").Value = 0.02336090049363864
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.751166666666666
$ws.Range("N2").Value = 17.2535
$ws.Range("O2").Value = 0.7405222614421495
$ws.Range("P2").Value = 0.7405222614421495
$ws.Range("Q2").Value = 0.545024645611111
$ws.Range("R2").Value = 4.9052218105
$ws.Range("S2").Value = 0.01729926686287431
$ws.Range("T2").Value = 0.01729926686287431
$ws.Range("G3").Value = 0.09476766666666665
$ws.Range("I3").Value = 0.02336090049363864
$ws.Range("J3").Value = 0.02336090049363864
$ws.Range("O3").Value = 0.07337387367415998
$ws.Range("P3").Value = 0.07337387367416
$ws.Range("Q3").Value = 0.05400319690388888
$ws.Range("R3").Value = 0.4860287721349999
$ws.Range("S3").Value = 0.001714079761734863
$ws.Range("T3").Value = 0.001714079761734863
$ws.Range("G4").Value = 0.09476766666666665
$ws.Range("I4").Value = 0.02336090049363864
$ws.Range("J4").Value = 0.02336090049363864
$ws.Range("M4").Value = 1.445350666666667
$ws.Range("N4").Value = 4.336052
$ws.Range("O4").Value = 0.1861038648836906
$ws.Range("P4").Value = 0.1861038648836906
$ws.Range("Q4").Value = 0.1369725101951111
$ws.Range("R4").Value = 1.232752591756
$ws.Range("S4").Value = 0.004347553869029466
$ws.Range("T4").Value = 0.004347553869029466
$ws.Range("H5").Value = 6.562189
$ws.Range("I5").Value = 0.5392086761288135
$ws.Range("J5").Value = 0.5392086761288134
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.751166666666666
$ws.Range("N5").Value = 17.2535
$ws.Range("O5").Value = 0.7405222614421495
$ws.Range("P5").Value = 0.7405222614421495
$ws.Range("Q5").Value = 12.58008087905556
$ws.Range("R5").Value = 113.2207279115
$ws.Range("S5").Value = 0.3992960282361366
$ws.Range("T5").Value = 0.3992960282361365
$ws.Range("H6").Value = 6.562189
$ws.Range("I6").Value = 0.5392086761288135
$ws.Range("J6").Value = 0.5392086761288134
$ws.Range("O6").Value = 0.07337387367415998
$ws.Range("P6").Value = 0.07337387367416
$ws.Range("Q6").Value = 1.246484154889445
$ws.Range("S6").Value = 0.03956382928628661
$ws.Range("T6").Value = 0.03956382928628661
$ws.Range("H7").Value = 6.562189
$ws.Range("I7").Value = 0.5392086761288135
$ws.Range("J7").Value = 0.5392086761288134
$ws.Range("M7").Value = 1.445350666666667
$ws.Range("N7").Value = 4.336052
$ws.Range("O7").Value = 0.1861038648836906
$ws.Range("P7").Value = 0.1861038648836906
$ws.Range("Q7").Value = 3.161554748647556
$ws.Range("R7").Value = 28.453992737828
$ws.Range("S7").Value = 0.1003488186063904
$ws.Range("T7").Value = 0.1003488186063904
$ws.Range("G8").Value = 1.774514666666667
$ws.Range("H8").Value = 5.323544
$ws.Range("I8").Value = 0.437430423377548
$ws.Range("J8").Value = 0.437430423377548
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.751166666666666
$ws.Range("N8").Value = 17.2535
$ws.Range("O8").Value = 0.7405222614421495
$ws.Range("P8").Value = 0.7405222614421495
$ws.Range("Q8").Value = 10.20552960044444
$ws.Range("R8").Value = 91.84976640399999
$ws.Range("S8").Value = 0.3239269663431387
$ws.Range("T8").Value = 0.3239269663431387
$ws.Range("G9").Value = 1.774514666666667
$ws.Range("H9").Value = 5.323544
$ws.Range("I9").Value = 0.437430423377548
$ws.Range("J9").Value = 0.437430423377548
$ws.Range("O9").Value = 0.07337387367415998
$ws.Range("P9").Value = 0.07337387367416
$ws.Range("Q9").Value = 1.011204225275556
$ws.Range("R9").Value = 9.10083802748
$ws.Range("S9").Value = 0.03209596462613852
$ws.Range("T9").Value = 0.03209596462613853
$ws.Range("G10").Value = 1.774514666666667
$ws.Range("H10").Value = 5.323544
$ws.Range("I10").Value = 0.437430423377548
$ws.Range("J10").Value = 0.437430423377548
$ws.Range("M10").Value = 1.445350666666667
$ws.Range("N10").Value = 4.336052
$ws.Range("O10").Value = 0.1861038648836906
$ws.Range("P10").Value = 0.1861038648836906
$ws.Range("Q10").Value = 2.564795956476445
$ws.Range("R10").Value = 23.083163608288
$ws.Range("S10").Value = 0.08140749240827076
$ws.Range("T10").Value = 0.08140749240827076
